$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.452.64"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.96%  "

$ws.Range("D3").Value = "1.825.62"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.74%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.20%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.67"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.24%  "

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.16%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5138"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -5.29%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3926"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +3.62%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07697"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.58%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.24%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.113"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.28%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.99"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +3.57%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.277"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.36%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.09%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.557"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.31%  "

$ws.Range("D16").Value = "1.824.78"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.59%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.34"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +5.15%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001079"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.84%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06638"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.10%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.99%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.06%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.27%  "

$ws.Range("D23").Value = "28.466.04"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.78%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.12"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.09%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.241"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +7.06%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.09"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.80%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.434"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +4.19%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.51%  "

$ws.Range("D29").Value = "2.035.44"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.56%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.75"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.74%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.129"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.15%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1096"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.72%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.642"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.25%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.643"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.65%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07190"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.34%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2230"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.23%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.943"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +6.06%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.87%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.155"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.96%  "

$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.28"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.56%  "

$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6238"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.83%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.190"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.99%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.002"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.09%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.393"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.83%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.44"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.23%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5898"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.29%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.706"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.63%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.33"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.03%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +3.59%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.48%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06932"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.99%  "
